$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- 1. Consolidate duplicate date-time number format (style 19 -> style 20) ---
# Copy the format from a pristine "style 20" cell (Sheet1!K3, numFmt 166 source,
# unaffected by any value edit below) onto every cell that currently carries the
# near-duplicate "style 19" (numFmt 164) format, so they end up sharing style 20.
$ws1.Range("K3").Copy()
$ws1.Range("D2").PasteSpecial(-4122)
$ws1.Range("J2").PasteSpecial(-4122)
$ws1.Range("D3:D4").PasteSpecial(-4122)
$ws1.Range("C5:D11").PasteSpecial(-4122)
$ws1.Range("D12:D14").PasteSpecial(-4122)
$ws1.Range("C15:D31").PasteSpecial(-4122)
$ws2.Range("D1:D3").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 2. Cell value / text updates ---
$ws1.Range("B3").Value = "Maximum Retry Error"
$ws1.Range("B4").Value = "Maximum Retry Error"

$ws1.Range("J2").Value = 45068.50023203704
$ws1.Range("K2").Value = 45068.52014737268
$ws1.Range("J3").Value = 45068.52014791666
$ws1.Range("J4").Value = 45068.5206337037
$ws1.Range("K4").Value = 45065.67844310185

Write-Host "Edit complete"
